# edit.ps1
# Applies the odds/statistics value updates described by the commit diff
# for "Jogos_da_Semana_FlashScore_2025-05-19.xlsx" (commit: "Atualizando o arquivo XLSX").
# All changes are plain numeric value overwrites on existing cells; no rows/
# columns are inserted or removed and no formatting/structure changes are made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 2.6
$ws.Range("H5").Value = 2.75
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 1.14
$ws.Range("K5").Value = 5.5
$ws.Range("L5").Value = 1.57
$ws.Range("M5").Value = 2.25
$ws.Range("N5").Value = 2.88
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 1.62
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 2.25
$ws.Range("S5").Value = 1.57
$ws.Range("T5").Value = 6
$ws.Range("X5").Value = 29
$ws.Range("Z5").Value = 5.5
$ws.Range("AC5").Value = 81
$ws.Range("AG5").Value = 13
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 51

# Row 10
$ws.Range("J10").Value = 1.05
$ws.Range("K10").Value = 11

# Row 17
$ws.Range("G17").Value = 2.3
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 2.8

# Row 20
$ws.Range("G20").Value = 2.52
$ws.Range("H20").Value = 3.7
$ws.Range("I20").Value = 2.42
$ws.Range("P20").Value = 1.26
$ws.Range("Q20").Value = 3.5
$ws.Range("T20").Value = 13.5
$ws.Range("U20").Value = 16.5
$ws.Range("V20").Value = 9.75
$ws.Range("W20").Value = 30
$ws.Range("X20").Value = 17.5
$ws.Range("Y20").Value = 19.5
$ws.Range("AE20").Value = 13.5
$ws.Range("AF20").Value = 16.5
$ws.Range("AG20").Value = 9.75
$ws.Range("AH20").Value = 29
$ws.Range("AI20").Value = 17

# Row 23
$ws.Range("I23").Value = 4.33
$ws.Range("K23").Value = 8
$ws.Range("L23").Value = 1.44
$ws.Range("M23").Value = 2.63
$ws.Range("R23").Value = 2.1
$ws.Range("S23").Value = 1.67
$ws.Range("AJ23").Value = 51

# Row 24
$ws.Range("G24").Value = 1.85
$ws.Range("I24").Value = 4.2
$ws.Range("J24").Value = 1.08
$ws.Range("K24").Value = 8
$ws.Range("R24").Value = 2
$ws.Range("S24").Value = 1.75
$ws.Range("T24").Value = 6
$ws.Range("U24").Value = 8
$ws.Range("Y24").Value = 34
$ws.Range("Z24").Value = 8
$ws.Range("AC24").Value = 67
$ws.Range("AD24").Value = 451
$ws.Range("AI24").Value = 41

# Row 25
$ws.Range("K25").Value = 19

# Row 26
$ws.Range("J26").Value = 1.02
$ws.Range("K26").Value = 21
$ws.Range("N26").Value = 1.36
$ws.Range("O26").Value = 3.1

# Row 28
$ws.Range("G28").Value = 1.26
$ws.Range("H28").Value = 5.1
$ws.Range("I28").Value = 8
$ws.Range("T28").Value = 7.2
$ws.Range("V28").Value = 7.6
$ws.Range("Z28").Value = 16
$ws.Range("AA28").Value = 9.25
$ws.Range("AB28").Value = 17.5
$ws.Range("AF28").Value = 45

# Row 31
$ws.Range("R31").Value = 1.67

# Row 32
$ws.Range("S32").Value = 1.57

# Row 33
$ws.Range("R33").Value = 1.75

# Row 34
$ws.Range("R34").Value = 1.67

# Row 35
$ws.Range("R35").Value = 1.67

# Row 36
$ws.Range("G36").Value = 2.63
$ws.Range("L36").Value = 1.3
$ws.Range("M36").Value = 3.4
$ws.Range("N36").Value = 2
$ws.Range("O36").Value = 1.8
$ws.Range("R36").Value = 1.8
$ws.Range("S36").Value = 1.91
$ws.Range("V36").Value = 10
$ws.Range("AA36").Value = 7
$ws.Range("AI36").Value = 19

# Row 43
$ws.Range("L43").Value = 1.2
$ws.Range("M43").Value = 4.33
$ws.Range("N43").Value = 1.67
$ws.Range("O43").Value = 2.15
$ws.Range("P43").Value = 1.3
$ws.Range("Q43").Value = 3.4
$ws.Range("T43").Value = 8.5
$ws.Range("U43").Value = 8.5
$ws.Range("AI43").Value = 41
$ws.Range("AJ43").Value = 41

# Row 48
$ws.Range("G48").Value = 1.35
$ws.Range("H48").Value = 5.2
$ws.Range("I48").Value = 6.9
$ws.Range("N48").Value = 1.39
$ws.Range("O48").Value = 2.77
$ws.Range("Q48").Value = 3.8
$ws.Range("W48").Value = 9.75
$ws.Range("X48").Value = 10
$ws.Range("AA48").Value = 11
$ws.Range("AE48").Value = 27
$ws.Range("AF48").Value = 55
$ws.Range("AG48").Value = 22
$ws.Range("AJ48").Value = 50

# Row 50
$ws.Range("G50").Value = 1.57
$ws.Range("H50").Value = 4.2
$ws.Range("I50").Value = 4.8
$ws.Range("K50").Value = 9.25
$ws.Range("L50").Value = 1.17
$ws.Range("M50").Value = 4.4
$ws.Range("N50").Value = 1.53
$ws.Range("O50").Value = 2.35
$ws.Range("P50").Value = 1.28
$ws.Range("Q50").Value = 3.35
$ws.Range("R50").Value = 1.6
$ws.Range("S50").Value = 2.2
$ws.Range("T50").Value = 9.25
$ws.Range("W50").Value = 12.5
$ws.Range("X50").Value = 11.25
$ws.Range("Y50").Value = 19.5
$ws.Range("Z50").Value = 9.25
$ws.Range("AA50").Value = 8.5
$ws.Range("AE50").Value = 17.5
$ws.Range("AF50").Value = 32
$ws.Range("AG50").Value = 15.5
$ws.Range("AH50").Value = 80
$ws.Range("AI50").Value = 40
$ws.Range("AJ50").Value = 37

# Row 51
$ws.Range("K51").Value = 9
